# Add team record (Wins / Losses / Ties) columns to the roster sheet.
# New columns: AD = Wins, AE = Losses, AF = Ties
# Header row (row 1) gets the same style as the existing header cells
# (bold, centered, bordered) by copying the format from AC1.
# Data rows (2-45) get the constant team record values: 77 wins, 85 losses, 0 ties.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 45

# Copy the existing header formatting (bold / border / alignment) onto the
# new header cells so they match the rest of row 1.
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Header labels.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Team record values for every data row.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 77
    $ws.Cells.Item($r, 31).Value = 85
    $ws.Cells.Item($r, 32).Value = 0
}
